$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 24 entry: "3. Danh mục sản phẩm (Phần 1)" as a hyperlink,
# mirroring the style of B22/B23 (existing list items in this section).
$url = "https://github.com/nguyentienminh07102004/product-management/commit/1fd9d226a6e73a845938d66ef4a661b05a626904"

# Create the hyperlink (this is how the cached/legacy "display" attribute on
# the <hyperlink> element ends up holding the raw commit URL), then replace
# the cell's visible text with the friendly Vietnamese title afterwards -
# exactly mirroring B22/B23 above it, whose "display" text is likewise
# stale relative to the text actually shown in the cell.
$ws.Hyperlinks.Add($ws.Range("B24"), $url, "", "", $url)
$ws.Range("B24").Value = "3. Danh mục sản phẩm (Phần 1)"

# Hyperlinks.Add stamps its own direct font formatting; re-apply the same
# "Hyperlink" cell style used by the rows above it (B22/B23) so B24 matches.
$ws.Range("B24").Style = $ws.Range("B23").Style

$ws.Range("B28").Select()
